# Auto-applied scheduled-runner price/profit refresh
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3332.2273
$ws.Range("I62").Value = 1599.9375
$ws.Range("J62").Value = 4322.107
$ws.Range("K62").Value = 1599.9375
$ws.Range("L62").Value = 4322.107
$ws.Range("M62").Value = -975.9375
$ws.Range("N62").Value = -5570.107
$ws.Range("H65").Value = 3332.2273
$ws.Range("I65").Value = 1599.9375
$ws.Range("J65").Value = 4322.107
$ws.Range("K65").Value = 7999.6875
$ws.Range("L65").Value = 21610.535
$ws.Range("M65").Value = -4879.6875
$ws.Range("N65").Value = -27850.535
$ws.Range("H98").Value = 2217.1316
$ws.Range("I98").Value = 2284.7778
$ws.Range("K98").Value = 2284.7778
$ws.Range("M98").Value = -786.7777999999998
$ws.Range("H122").Value = 2217.1316
$ws.Range("I122").Value = 2284.7778
$ws.Range("K122").Value = 6854.3334
$ws.Range("M122").Value = -4404.3334
$ws.Range("H129").Value = 913.1739
$ws.Range("I129").Value = 980
$ws.Range("J129").Value = 910.13635
$ws.Range("K129").Value = 2940
$ws.Range("L129").Value = 2730.40905
$ws.Range("N129").Value = -12730.40905
$ws.Range("M129").Value = 2060
$ws.Range("H132").Value = 2165.0747
$ws.Range("I132").Value = 2207.541
$ws.Range("J132").Value = 1733.3334
$ws.Range("K132").Value = 6622.623000000001
$ws.Range("L132").Value = 5200.0002
$ws.Range("M132").Value = -4092.623000000001
$ws.Range("N132").Value = -10260.0002
$ws.Range("H138").Value = 1332.54
$ws.Range("I138").Value = 667.9286
$ws.Range("J138").Value = 2883.3
$ws.Range("K138").Value = 2003.7858
$ws.Range("L138").Value = 8649.900000000001
$ws.Range("M138").Value = 3136.2142
$ws.Range("N138").Value = -18929.9

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2464.74
$ws.Range("I32").Value = 1947.9375
$ws.Range("J32").Value = 4531.95
$ws.Range("K32").Value = 1947.9375
$ws.Range("L32").Value = 4531.95
$ws.Range("M32").Value = -1660.9375
$ws.Range("N32").Value = -5105.95
$ws.Range("H45").Value = 1046.7778
$ws.Range("I45").Value = 721.4
$ws.Range("J45").Value = 1453.5
$ws.Range("K45").Value = 721.4
$ws.Range("L45").Value = 1453.5
$ws.Range("M45").Value = -344.4
$ws.Range("N45").Value = -2207.5
$ws.Range("H61").Value = 2179700.8
$ws.Range("I61").Value = 2850034.5
$ws.Range("J61").Value = 1116.3334
$ws.Range("K61").Value = 2850034.5
$ws.Range("L61").Value = 1116.3334
$ws.Range("M61").Value = -2849822.5
$ws.Range("N61").Value = -1540.3334
$ws.Range("H110").Value = 1379.9546
$ws.Range("I110").Value = 1139.5333
$ws.Range("J110").Value = 1895.1428
$ws.Range("K110").Value = 1139.5333
$ws.Range("L110").Value = 1895.1428
$ws.Range("M110").Value = 905.4666999999999
$ws.Range("N110").Value = -5985.1428
$ws.Range("H122").Value = 1416.9269
$ws.Range("I122").Value = 1355.4073
$ws.Range("J122").Value = 1535.5714
$ws.Range("K122").Value = 4066.2219
$ws.Range("L122").Value = 4606.7142
$ws.Range("M122").Value = -1616.2219
$ws.Range("N122").Value = -9506.7142
$ws.Range("H132").Value = 6910488
$ws.Range("I132").Value = 8928924
$ws.Range("J132").Value = 78858.08
$ws.Range("K132").Value = 26786772
$ws.Range("L132").Value = 236574.24
$ws.Range("M132").Value = -26784242
$ws.Range("N132").Value = -241634.24
$ws.Range("H133").Value = 33753.668
$ws.Range("J133").Value = 33753.668
$ws.Range("L133").Value = 33753.668
$ws.Range("N133").Value = -38813.668
$ws.Range("H136").Value = 2179700.8
$ws.Range("I136").Value = 2850034.5
$ws.Range("J136").Value = 1116.3334
$ws.Range("K136").Value = 8550103.5
$ws.Range("L136").Value = 3349.0002
$ws.Range("M136").Value = -8547553.5
$ws.Range("N136").Value = -8449.0002
$ws.Range("H139").Value = 35560.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 35560.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 35560.75
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -45840.75
$ws.Range("H140").Value = 40241.4
$ws.Range("J140").Value = 40241.4
$ws.Range("L140").Value = 40241.4
$ws.Range("N140").Value = -50601.4

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 40196.47
$ws.Range("J135").Value = 40196.47
$ws.Range("L135").Value = 40196.47
$ws.Range("N135").Value = -50336.47
$ws.Range("H140").Value = 70772
$ws.Range("J140").Value = 70772
$ws.Range("L140").Value = 70772
$ws.Range("N140").Value = -81132

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1407.2354
$ws.Range("I134").Value = 1676.9131
$ws.Range("J134").Value = 843.36365
$ws.Range("K134").Value = 5030.7393
$ws.Range("L134").Value = 2530.09095
$ws.Range("M134").Value = -2495.7393
$ws.Range("N134").Value = -7600.09095
$ws.Range("H140").Value = 39390
$ws.Range("I140").Value = 30000
$ws.Range("J140").Value = 48780
$ws.Range("K140").Value = 30000
$ws.Range("L140").Value = 48780
$ws.Range("M140").Value = -24820
$ws.Range("N140").Value = -59140

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8929391
$ws.Range("I5").Value = 426.6842
$ws.Range("J5").Value = 13514536
$ws.Range("K5").Value = 1280.0526
$ws.Range("L5").Value = 40543608
$ws.Range("M5").Value = -1168.0526
$ws.Range("N5").Value = -40543832
$ws.Range("H122").Value = 13228961
$ws.Range("I122").Value = 22222958
$ws.Range("J122").Value = 1986464.5
$ws.Range("K122").Value = 200006622
$ws.Range("L122").Value = 17878180.5
$ws.Range("M122").Value = -200004172
$ws.Range("N122").Value = -17883080.5
$ws.Range("H131").Value = 5983.476
$ws.Range("J131").Value = 4790.231
$ws.Range("L131").Value = 14370.693
$ws.Range("N131").Value = -24450.693
$ws.Range("H135").Value = 8929391
$ws.Range("I135").Value = 426.6842
$ws.Range("J135").Value = 13514536
$ws.Range("K135").Value = 3840.1578
$ws.Range("L135").Value = 121630824
$ws.Range("M135").Value = -1305.1578
$ws.Range("N135").Value = -121635894

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 284.2857
$ws.Range("J22").Value = 297.75
$ws.Range("L22").Value = 297.75
$ws.Range("N22").Value = -887.75
$ws.Range("H27").Value = 284.2857
$ws.Range("J27").Value = 297.75
$ws.Range("L27").Value = 297.75
$ws.Range("N27").Value = -511.75
$ws.Range("H61").Value = 1177.9231
$ws.Range("I61").Value = 1243
$ws.Range("J61").Value = 820
$ws.Range("K61").Value = 1243
$ws.Range("L61").Value = 820
$ws.Range("M61").Value = -1041
$ws.Range("N61").Value = -1224
$ws.Range("H113").Value = 1177.9231
$ws.Range("I113").Value = 1243
$ws.Range("J113").Value = 820
$ws.Range("K113").Value = 1243
$ws.Range("L113").Value = 820
$ws.Range("M113").Value = 927
$ws.Range("N113").Value = -5160

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 22333.334
$ws.Range("J46").Value = 22333.334
$ws.Range("L46").Value = 22333.334
$ws.Range("M46").Value = -22795.334
$ws.Range("H113").Value = 438.56
$ws.Range("I113").Value = 399.06668
$ws.Range("J113").Value = 497.8
$ws.Range("K113").Value = 1197.20004
$ws.Range("L113").Value = 1493.4
$ws.Range("M113").Value = 972.7999599999998
$ws.Range("N113").Value = -5833.4
$ws.Range("H132").Value = 6227941
$ws.Range("I132").Value = 7265524
$ws.Range("J132").Value = 2444.3333
$ws.Range("K132").Value = 21796572
$ws.Range("L132").Value = 7332.999899999999
$ws.Range("M132").Value = -21794042
$ws.Range("N132").Value = -12392.9999
$ws.Range("H134").Value = 22333.334
$ws.Range("J134").Value = 22333.334
$ws.Range("L134").Value = 67000.00199999999
$ws.Range("N134").Value = -72070.00199999999
$ws.Range("H136").Value = 2752077.5
$ws.Range("I136").Value = 6383
$ws.Range("J136").Value = 8404978
$ws.Range("K136").Value = 19149
$ws.Range("L136").Value = 25214934
$ws.Range("M136").Value = -16599
$ws.Range("N136").Value = -25220034
